$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = 6
$ws.Range("B11").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C11").Value = "Metropolitana"
$ws.Range("D11").Value = 44881
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100108
$ws.Range("H11").Value = "Tropicales y subtropicales"
$ws.Range("I11").Value = 100108007
$ws.Range("J11").Value = "Coco"
$ws.Range("K11").Value = "Sin especificar"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 80
$ws.Range("N11").Value = 28000
$ws.Range("O11").Value = 30000
$ws.Range("P11").Value = 29000
$ws.Range("Q11").Value = '$/malla 20 unidades'
$ws.Range("R11").Value = "Perú"
$ws.Range("S11").Value = 1450
$ws.Range("T11").Value = 20
